$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 7002.857
$ws.Range("J45").Value = 250
$ws.Range("L45").Value = 750
$ws.Range("N45").Value = -1134

$ws.Range("H48").Value = 399.5
$ws.Range("J48").Value = 499
$ws.Range("L48").Value = 1497
$ws.Range("N48").Value = -2081

$ws.Range("H56").Value = 399.5
$ws.Range("J56").Value = 499
$ws.Range("L56").Value = 1497
$ws.Range("N56").Value = -2565

$ws.Range("H59").Value = 625
$ws.Range("I59").Value = 1000
$ws.Range("J59").Value = 250
$ws.Range("K59").Value = 3000
$ws.Range("L59").Value = 750
$ws.Range("M59").Value = -2443
$ws.Range("N59").Value = -1864

$ws.Range("H116").Value = 3593.2693
$ws.Range("I116").Value = 3672.476
$ws.Range("J116").Value = 3260.6
$ws.Range("K116").Value = 3672.476
$ws.Range("L116").Value = 3260.6
$ws.Range("M116").Value = -230.4760000000001
$ws.Range("N116").Value = -10144.6

$ws.Range("H125").Value = 520.9
$ws.Range("I125").Value = 429.66666
$ws.Range("J125").Value = 657.75
$ws.Range("K125").Value = 3866.99994
$ws.Range("L125").Value = 5919.75
$ws.Range("M125").Value = -1406.99994
$ws.Range("N125").Value = -10839.75

$ws.Range("H131").Value = 2362.8125
$ws.Range("I131").Value = 2362.8125
$ws.Range("K131").Value = 7088.4375
$ws.Range("M131").Value = -2048.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5161.355
$ws.Range("I32").Value = 3550.8076
$ws.Range("K32").Value = 3550.8076
$ws.Range("M32").Value = -3263.8076

$ws.Range("H74").Value = 2739.7917
$ws.Range("I74").Value = 2438.3
$ws.Range("K74").Value = 2438.3
$ws.Range("M74").Value = -1564.3

$ws.Range("H77").Value = 2739.7917
$ws.Range("I77").Value = 2438.3
$ws.Range("K77").Value = 12191.5
$ws.Range("M77").Value = -7823.5

$ws.Range("H112").Value = 43000
$ws.Range("J112").Value = 43000
$ws.Range("L112").Value = 43000
$ws.Range("N112").Value = -45954

$ws.Range("H122").Value = 5560610.5
$ws.Range("I122").Value = 6671776
$ws.Range("J122").Value = 4783.1665
$ws.Range("K122").Value = 20015328
$ws.Range("L122").Value = 14349.4995
$ws.Range("M122").Value = -20012878
$ws.Range("N122").Value = -19249.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3322.1333
$ws.Range("I86").Value = 3068.4443
$ws.Range("J86").Value = 3702.6667
$ws.Range("K86").Value = 3068.4443
$ws.Range("L86").Value = 3702.6667
$ws.Range("M86").Value = -1945.4443
$ws.Range("N86").Value = -5948.6667

$ws.Range("H89").Value = 3322.1333
$ws.Range("I89").Value = 3068.4443
$ws.Range("J89").Value = 3702.6667
$ws.Range("K89").Value = 15342.2215
$ws.Range("L89").Value = 18513.3335
$ws.Range("M89").Value = -9726.2215
$ws.Range("N89").Value = -29745.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2384
$ws.Range("I31").Value = 2289.889
$ws.Range("J31").Value = 2626
$ws.Range("K31").Value = 2289.889
$ws.Range("L31").Value = 2626
$ws.Range("M31").Value = -1994.889
$ws.Range("N31").Value = -3216

$ws.Range("H34").Value = 2384
$ws.Range("I34").Value = 2289.889
$ws.Range("J34").Value = 2626
$ws.Range("K34").Value = 2289.889
$ws.Range("L34").Value = 2626
$ws.Range("M34").Value = -2087.889
$ws.Range("N34").Value = -3030

$ws.Range("H58").Value = 3741.88
$ws.Range("I58").Value = 3820.3635
$ws.Range("K58").Value = 3820.3635
$ws.Range("M58").Value = -3617.3635

$ws.Range("H132").Value = 2081.7727
$ws.Range("I132").Value = 1792.6154
$ws.Range("J132").Value = 2499.4443
$ws.Range("K132").Value = 5377.8462
$ws.Range("L132").Value = 7498.3329
$ws.Range("M132").Value = -2847.8462
$ws.Range("N132").Value = -12558.3329

$ws.Range("H134").Value = 1184.875
$ws.Range("I134").Value = 1184.875
$ws.Range("K134").Value = 3554.625
$ws.Range("M134").Value = -1019.625

$ws.Range("H136").Value = 3741.88
$ws.Range("I136").Value = 3820.3635
$ws.Range("K136").Value = 11461.0905
$ws.Range("M136").Value = -8911.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 17499.5
$ws.Range("I68").Value = 17499.5
$ws.Range("K68").Value = 52498.5
$ws.Range("M68").Value = -51687.5

$ws.Range("H71").Value = 17499.5
$ws.Range("I71").Value = 17499.5
$ws.Range("K71").Value = 157495.5
$ws.Range("M71").Value = -153439.5

$ws.Range("H117").Value = 45923.332
$ws.Range("I117").Value = 399
$ws.Range("J117").Value = 47902.652
$ws.Range("K117").Value = 1197
$ws.Range("L117").Value = 143707.956
$ws.Range("M117").Value = 2245
$ws.Range("N117").Value = -150591.956

$ws.Range("H122").Value = 1178.3572
$ws.Range("I122").Value = 1416.3334
$ws.Range("J122").Value = 750
$ws.Range("K122").Value = 12747.0006
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -10297.0006
$ws.Range("N122").Value = -11650

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1492
$ws.Range("I97").Value = 1931.5
$ws.Range("J97").Value = 1199
$ws.Range("K97").Value = 1931.5
$ws.Range("L97").Value = 1199
$ws.Range("M97").Value = -1435.5
$ws.Range("N97").Value = -2191

$ws.Range("H122").Value = 47621108
$ws.Range("I122").Value = 55556680
$ws.Range("K122").Value = 166670040
$ws.Range("M122").Value = -166667590

$ws.Range("H123").Value = 93999
$ws.Range("J123").Value = 93999
$ws.Range("L123").Value = 93999
$ws.Range("N123").Value = -98899

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 40084
$ws.Range("J51").Value = 40084
$ws.Range("L51").Value = 40084
$ws.Range("N51").Value = -41040

$ws.Range("H110").Value = 38000
$ws.Range("J110").Value = 38000
$ws.Range("L110").Value = 38000
$ws.Range("N110").Value = -46180

$ws.Range("H132").Value = 5361.4053
$ws.Range("J132").Value = 5234.5
$ws.Range("L132").Value = 15703.5
$ws.Range("N132").Value = -20763.5

$ws.Range("H136").Value = 5374.825
$ws.Range("I136").Value = 4249.7666
$ws.Range("J136").Value = 8750
$ws.Range("K136").Value = 12749.2998
$ws.Range("L136").Value = 26250
$ws.Range("M136").Value = -10199.2998
$ws.Range("N136").Value = -31350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 7500
$ws.Range("I51").Value = 7500
$ws.Range("K51").Value = 7500
$ws.Range("M51").Value = -6990

$ws.Range("H96").Value = 12347253
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H126").Value = 2386.9
$ws.Range("I126").Value = 2158.7144
$ws.Range("J126").Value = 2919.3333
$ws.Range("K126").Value = 6476.1432
$ws.Range("L126").Value = 8757.999899999999
$ws.Range("M126").Value = -4006.1432
$ws.Range("N126").Value = -13697.9999

$ws.Range("H132").Value = 1643.0869
$ws.Range("I132").Value = 1608.6818
$ws.Range("K132").Value = 4826.0454
$ws.Range("M132").Value = -2296.0454

$ws.Range("H136").Value = 8031
$ws.Range("I136").Value = 6705.5835
$ws.Range("K136").Value = 20116.7505
$ws.Range("M136").Value = -17566.7505
